$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.555934
$ws.Cells.Item(2, 8).Value = 1.667802
$ws.Cells.Item(2, 9).Value = 0.005745252779589096
$ws.Cells.Item(2, 10).Value = 0.005745252779589094
$ws.Cells.Item(2, 13).Value = 133.7780026666667
$ws.Cells.Item(2, 14).Value = 401.334008
$ws.Cells.Item(2, 15).Value = 0.50863533211804
$ws.Cells.Item(2, 16).Value = 0.5086353321180399
$ws.Cells.Item(2, 17).Value = 74.37174013449068
$ws.Cells.Item(2, 18).Value = 669.3456612104161
$ws.Cells.Item(2, 19).Value = 0.002922238555648392
$ws.Cells.Item(2, 20).Value = 0.002922238555648391
$ws.Cells.Item(3, 7).Value = 0.555934
$ws.Cells.Item(3, 8).Value = 1.667802
$ws.Cells.Item(3, 9).Value = 0.005745252779589096
$ws.Cells.Item(3, 10).Value = 0.005745252779589094
$ws.Cells.Item(3, 15).Value = 0.1993888292903622
$ws.Cells.Item(3, 16).Value = 0.1993888292903622
$ws.Cells.Item(3, 17).Value = 29.154274705914
$ws.Cells.Item(3, 18).Value = 262.388472353226
$ws.Cells.Item(3, 19).Value = 0.001145539225699469
$ws.Cells.Item(3, 20).Value = 0.001145539225699469
$ws.Cells.Item(4, 7).Value = 0.555934
$ws.Cells.Item(4, 8).Value = 1.667802
$ws.Cells.Item(4, 9).Value = 0.005745252779589096
$ws.Cells.Item(4, 10).Value = 0.005745252779589094
$ws.Cells.Item(4, 13).Value = 21.197691
$ws.Cells.Item(4, 14).Value = 63.593073
$ws.Cells.Item(4, 15).Value = 0.08059542216956049
$ws.Cells.Item(4, 16).Value = 0.08059542216956046
$ws.Cells.Item(4, 17).Value = 11.784517148394
$ws.Cells.Item(4, 18).Value = 106.060654335546
$ws.Cells.Item(4, 19).Value = 0.000463041073241824
$ws.Cells.Item(4, 20).Value = 0.0004630410732418237
$ws.Cells.Item(5, 7).Value = 0.555934
$ws.Cells.Item(5, 8).Value = 1.667802
$ws.Cells.Item(5, 9).Value = 0.005745252779589096
$ws.Cells.Item(5, 10).Value = 0.005745252779589094
$ws.Cells.Item(5, 13).Value = 55.59592133333333
$ws.Cells.Item(5, 14).Value = 166.787764
$ws.Cells.Item(5, 15).Value = 0.2113804164220374
$ws.Cells.Item(5, 16).Value = 0.2113804164220373
$ws.Cells.Item(5, 17).Value = 30.90766293052533
$ws.Cells.Item(5, 18).Value = 278.168966374728
$ws.Cells.Item(5, 19).Value = 0.001214433924999411
$ws.Cells.Item(5, 20).Value = 0.00121443392499941
$ws.Cells.Item(6, 9).Value = 0.823525905561055
$ws.Cells.Item(6, 10).Value = 0.823525905561055
$ws.Cells.Item(6, 13).Value = 133.7780026666667
$ws.Cells.Item(6, 14).Value = 401.334008
$ws.Cells.Item(6, 15).Value = 0.50863533211804
$ws.Cells.Item(6, 16).Value = 0.5086353321180399
$ws.Cells.Item(6, 17).Value = 10660.46299302923
$ws.Cells.Item(6, 18).Value = 95944.16693726311
$ws.Cells.Item(6, 19).Value = 0.4188743724828569
$ws.Cells.Item(6, 20).Value = 0.4188743724828567
$ws.Cells.Item(7, 9).Value = 0.823525905561055
$ws.Cells.Item(7, 10).Value = 0.823525905561055
$ws.Cells.Item(7, 15).Value = 0.1993888292903622
$ws.Cells.Item(7, 16).Value = 0.1993888292903622
$ws.Cells.Item(7, 19).Value = 0.1642018662001042
$ws.Cells.Item(7, 20).Value = 0.1642018662001041
$ws.Cells.Item(8, 9).Value = 0.823525905561055
$ws.Cells.Item(8, 10).Value = 0.823525905561055
$ws.Cells.Item(8, 13).Value = 21.197691
$ws.Cells.Item(8, 14).Value = 63.593073
$ws.Cells.Item(8, 15).Value = 0.08059542216956049
$ws.Cells.Item(8, 16).Value = 0.08059542216956046
$ws.Cells.Item(8, 17).Value = 1689.195502538889
$ws.Cells.Item(8, 18).Value = 15202.75952285
$ws.Cells.Item(8, 19).Value = 0.06637241802626283
$ws.Cells.Item(8, 20).Value = 0.0663724180262628
$ws.Cells.Item(9, 9).Value = 0.823525905561055
$ws.Cells.Item(9, 10).Value = 0.823525905561055
$ws.Cells.Item(9, 13).Value = 55.59592133333333
$ws.Cells.Item(9, 14).Value = 166.787764
$ws.Cells.Item(9, 15).Value = 0.2113804164220374
$ws.Cells.Item(9, 16).Value = 0.2113804164220373
$ws.Cells.Item(9, 17).Value = 4430.311786117296
$ws.Cells.Item(9, 18).Value = 39872.80607505567
$ws.Cells.Item(9, 19).Value = 0.1740772488518312
$ws.Cells.Item(9, 20).Value = 0.1740772488518312
$ws.Cells.Item(10, 7).Value = 0.3446996666666666
$ws.Cells.Item(10, 8).Value = 1.034099
$ws.Cells.Item(10, 9).Value = 0.003562269474506148
$ws.Cells.Item(10, 10).Value = 0.003562269474506148
$ws.Cells.Item(10, 13).Value = 133.7780026666667
$ws.Cells.Item(10, 14).Value = 401.334008
$ws.Cells.Item(10, 15).Value = 0.50863533211804
$ws.Cells.Item(10, 16).Value = 0.5086353321180399
$ws.Cells.Item(10, 17).Value = 46.11323292653244
$ws.Cells.Item(10, 18).Value = 415.019096338792
$ws.Cells.Item(10, 19).Value = 0.00181189611725939
$ws.Cells.Item(10, 20).Value = 0.00181189611725939
$ws.Cells.Item(11, 7).Value = 0.3446996666666666
$ws.Cells.Item(11, 8).Value = 1.034099
$ws.Cells.Item(11, 9).Value = 0.003562269474506148
$ws.Cells.Item(11, 10).Value = 0.003562269474506148
$ws.Cells.Item(11, 15).Value = 0.1993888292903622
$ws.Cells.Item(11, 16).Value = 0.1993888292903622
$ws.Cells.Item(11, 17).Value = 18.076729923043
$ws.Cells.Item(11, 18).Value = 162.690569307387
$ws.Cells.Item(11, 19).Value = 0.0007102767401385747
$ws.Cells.Item(11, 20).Value = 0.0007102767401385746
$ws.Cells.Item(12, 7).Value = 0.3446996666666666
$ws.Cells.Item(12, 8).Value = 1.034099
$ws.Cells.Item(12, 9).Value = 0.003562269474506148
$ws.Cells.Item(12, 10).Value = 0.003562269474506148
$ws.Cells.Item(12, 13).Value = 21.197691
$ws.Cells.Item(12, 14).Value = 63.593073
$ws.Cells.Item(12, 15).Value = 0.08059542216956049
$ws.Cells.Item(12, 16).Value = 0.08059542216956046
$ws.Cells.Item(12, 17).Value = 7.306837021803
$ws.Cells.Item(12, 18).Value = 65.76153319622699
$ws.Cells.Item(12, 19).Value = 0.0002871026121795614
$ws.Cells.Item(12, 20).Value = 0.0002871026121795613
$ws.Cells.Item(13, 7).Value = 0.3446996666666666
$ws.Cells.Item(13, 8).Value = 1.034099
$ws.Cells.Item(13, 9).Value = 0.003562269474506148
$ws.Cells.Item(13, 10).Value = 0.003562269474506148
$ws.Cells.Item(13, 13).Value = 55.59592133333333
$ws.Cells.Item(13, 14).Value = 166.787764
$ws.Cells.Item(13, 15).Value = 0.2113804164220374
$ws.Cells.Item(13, 16).Value = 0.2113804164220373
$ws.Cells.Item(13, 17).Value = 19.16389555162622
$ws.Cells.Item(13, 18).Value = 172.475059964636
$ws.Cells.Item(13, 19).Value = 0.0007529940049286218
$ws.Cells.Item(13, 20).Value = 0.0007529940049286215
$ws.Cells.Item(14, 7).Value = 16.17571666666667
$ws.Cells.Item(14, 8).Value = 48.52715
$ws.Cells.Item(14, 9).Value = 0.1671665721848498
$ws.Cells.Item(14, 10).Value = 0.1671665721848498
$ws.Cells.Item(14, 13).Value = 133.7780026666667
$ws.Cells.Item(14, 14).Value = 401.334008
$ws.Cells.Item(14, 15).Value = 0.50863533211804
$ws.Cells.Item(14, 16).Value = 0.5086353321180399
$ws.Cells.Item(14, 17).Value = 2163.955067368578
$ws.Cells.Item(14, 18).Value = 19475.5956063172
$ws.Cells.Item(14, 19).Value = 0.08502682496227541
$ws.Cells.Item(14, 20).Value = 0.08502682496227537
$ws.Cells.Item(15, 7).Value = 16.17571666666667
$ws.Cells.Item(15, 8).Value = 48.52715
$ws.Cells.Item(15, 9).Value = 0.1671665721848498
$ws.Cells.Item(15, 10).Value = 0.1671665721848498
$ws.Cells.Item(15, 15).Value = 0.1993888292903622
$ws.Cells.Item(15, 16).Value = 0.1993888292903622
$ws.Cells.Item(15, 17).Value = 848.28646433755
$ws.Cells.Item(15, 18).Value = 7634.578179037951
$ws.Cells.Item(15, 19).Value = 0.03333114712442004
$ws.Cells.Item(15, 20).Value = 0.03333114712442003
$ws.Cells.Item(16, 7).Value = 16.17571666666667
$ws.Cells.Item(16, 8).Value = 48.52715
$ws.Cells.Item(16, 9).Value = 0.1671665721848498
$ws.Cells.Item(16, 10).Value = 0.1671665721848498
$ws.Cells.Item(16, 13).Value = 21.197691
$ws.Cells.Item(16, 14).Value = 63.593073
$ws.Cells.Item(16, 15).Value = 0.08059542216956049
$ws.Cells.Item(16, 16).Value = 0.08059542216956046
$ws.Cells.Item(16, 17).Value = 342.88784360355
$ws.Cells.Item(16, 18).Value = 3085.99059243195
$ws.Cells.Item(16, 19).Value = 0.01347286045787628
$ws.Cells.Item(16, 20).Value = 0.01347286045787627
$ws.Cells.Item(17, 7).Value = 16.17571666666667
$ws.Cells.Item(17, 8).Value = 48.52715
$ws.Cells.Item(17, 9).Value = 0.1671665721848498
$ws.Cells.Item(17, 10).Value = 0.1671665721848498
$ws.Cells.Item(17, 13).Value = 55.59592133333333
$ws.Cells.Item(17, 14).Value = 166.787764
$ws.Cells.Item(17, 15).Value = 0.2113804164220374
$ws.Cells.Item(17, 16).Value = 0.2113804164220373
$ws.Cells.Item(17, 17).Value = 899.3038713102889
$ws.Cells.Item(17, 18).Value = 8093.734841792599
$ws.Cells.Item(17, 19).Value = 0.03533573964027813
$ws.Cells.Item(17, 20).Value = 0.03533573964027811
